# Rename the worksheet from "Sheet1" to "Estimates".
# (Excel auto-updates the _xlnm._FilterDatabase defined names that
# reference the sheet name when the sheet is renamed.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Estimates"

# Cell A87 currently duplicates A86's label "Min (P=95%)" - it should
# actually read "Max (P=95%)" (a new shared string gets added).
$ws.Range("A87").Value = "Max (P=95%)"
